$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DW")

$ws.Range("A27").Value = "subtree problem"
$ws.Range("B27").Value = "Tree"
$ws.Range("C27").Value = "Easy"
$ws.Range("D27").Value = "https://leetcode.com/problems/subtree-of-another-tree/"
$ws.Range("E27").Value = "Use Sametree function and check for edge case cases like null and all and recursively call itself too "

$ws.Range("A28").Value = "Reverse Word in a string"
$ws.Range("B28").Value = "String"
$ws.Range("C28").Value = "Medium"
$ws.Range("D28").Value = "https://leetcode.com/problems/reverse-words-in-a-string/"
$ws.Range("E28").Value = "Reverse whole word then iterate and make words  and add "

$ws.Range("A28").Select()
